$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new value. Using NumberFormat "@" (Text) before
# assignment and clearing it afterwards keeps Excel from re-interpreting
# numeric-looking strings (prices, percentages) as actual numbers, so values
# like "0.5090" or "28.028.59" are preserved exactly as literal text.
$updates = @(
    @{ Cell = 'D2'; Value = '28.028.59' }
    @{ Cell = 'E2'; Value = '  -0.38%  ' }
    @{ Cell = 'D3'; Value = '1.862.41' }
    @{ Cell = 'E3'; Value = '  -0.70%  ' }
    @{ Cell = 'E4'; Value = '  +0.50%  ' }
    @{ Cell = 'D5'; Value = '312.38' }
    @{ Cell = 'E5'; Value = '  -0.14%  ' }
    @{ Cell = 'E6'; Value = '  +0.43%  ' }
    @{ Cell = 'D7'; Value = '0.5090' }
    @{ Cell = 'E7'; Value = '  +0.91%  ' }
    @{ Cell = 'D8'; Value = '0.3852' }
    @{ Cell = 'E8'; Value = '  +0.18%  ' }
    @{ Cell = 'D9'; Value = '0.08272' }
    @{ Cell = 'E9'; Value = '  -8.29%  ' }
    @{ Cell = 'D10'; Value = '1.115' }
    @{ Cell = 'E10'; Value = '  -0.56%  ' }
    @{ Cell = 'B11'; Value = 'OKB' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D11'; Value = '41.44' }
    @{ Cell = 'E11'; Value = '  -0.66%  ' }
    @{ Cell = 'B12'; Value = 'Polkadot' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D12'; Value = '6.225' }
    @{ Cell = 'E12'; Value = '  -2.32%  ' }
    @{ Cell = 'B13'; Value = 'Solana' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' }
    @{ Cell = 'D13'; Value = '20.56' }
    @{ Cell = 'E13'; Value = '  -1.04%  ' }
    @{ Cell = 'B14'; Value = 'WrappedEther' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D14'; Value = '1.862.54' }
    @{ Cell = 'E14'; Value = '  -0.63%  ' }
    @{ Cell = 'B15'; Value = 'Chainlink' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = 'D15'; Value = '7.232' }
    @{ Cell = 'E15'; Value = '  -0.58%  ' }
    @{ Cell = 'B16'; Value = 'BinanceUSD' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' }
    @{ Cell = 'D16'; Value = '1.004' }
    @{ Cell = 'E16'; Value = '  +0.48%  ' }
    @{ Cell = 'B17'; Value = 'ShibaInu' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D17'; Value = '0.00001097' }
    @{ Cell = 'E17'; Value = '  -0.99%  ' }
    @{ Cell = 'B18'; Value = 'Litecoin' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D18'; Value = '90.81' }
    @{ Cell = 'E18'; Value = '  -0.56%  ' }
    @{ Cell = 'B19'; Value = 'TRON' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' }
    @{ Cell = 'D19'; Value = '0.06642' }
    @{ Cell = 'E19'; Value = '  -0.15%  ' }
    @{ Cell = 'B20'; Value = 'Avalanche' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'D20'; Value = '17.74' }
    @{ Cell = 'E20'; Value = '  -2.65%  ' }
    @{ Cell = 'B21'; Value = 'Dai' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D21'; Value = '1.003' }
    @{ Cell = 'E21'; Value = '  +0.34%  ' }
    @{ Cell = 'B22'; Value = 'Uniswap' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D22'; Value = '6.039' }
    @{ Cell = 'E22'; Value = '  -1.66%  ' }
    @{ Cell = 'B23'; Value = 'WrappedBTC' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = 'D23'; Value = '28.047.21' }
    @{ Cell = 'E23'; Value = '  -0.37%  ' }
    @{ Cell = 'B24'; Value = 'Cosmos' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D24'; Value = '11.11' }
    @{ Cell = 'E24'; Value = '  -3.21%  ' }
    @{ Cell = 'B25'; Value = 'Toncoin' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D25'; Value = '2.230' }
    @{ Cell = 'E25'; Value = '  -0.96%  ' }
    @{ Cell = 'B26'; Value = 'LidoDAOToken' }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D26'; Value = '2.538' }
    @{ Cell = 'E26'; Value = '  -0.23%  ' }
    @{ Cell = 'B27'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D27'; Value = '2.074.86' }
    @{ Cell = 'E27'; Value = '  -0.53%  ' }
    @{ Cell = 'B28'; Value = 'Monero' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D28'; Value = '157.73' }
    @{ Cell = 'E28'; Value = '  +0.45%  ' }
    @{ Cell = 'B29'; Value = 'EthereumClassic' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D29'; Value = '20.53' }
    @{ Cell = 'E29'; Value = '  -1.51%  ' }
    @{ Cell = 'B30'; Value = 'BitcoinCash' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D30'; Value = '124.86' }
    @{ Cell = 'E30'; Value = '  -1.62%  ' }
    @{ Cell = 'B31'; Value = 'Stellar' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Cell = 'D31'; Value = '0.1059' }
    @{ Cell = 'E31'; Value = '  -0.54%  ' }
    @{ Cell = 'B32'; Value = 'ImmutableX' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D32'; Value = '1.037' }
    @{ Cell = 'E32'; Value = '  -2.57%  ' }
    @{ Cell = 'B33'; Value = 'Filecoin' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D33'; Value = '5.902' }
    @{ Cell = 'E33'; Value = '  +5.11%  ' }
    @{ Cell = 'B34'; Value = 'HuobiToken' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D34'; Value = '3.592' }
    @{ Cell = 'E34'; Value = '  +0.02%  ' }
    @{ Cell = 'B35'; Value = 'FraxShare' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D35'; Value = '9.413' }
    @{ Cell = 'E35'; Value = '  -0.38%  ' }
    @{ Cell = 'B36'; Value = 'Hedera' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D36'; Value = '0.06555' }
    @{ Cell = 'E36'; Value = '  -0.52%  ' }
    @{ Cell = 'B37'; Value = 'VeChain' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D37'; Value = '0.02420' }
    @{ Cell = 'E37'; Value = '  +0.50%  ' }
    @{ Cell = 'B38'; Value = 'Algorand' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D38'; Value = '0.2177' }
    @{ Cell = 'E38'; Value = '  -0.83%  ' }
    @{ Cell = 'B39'; Value = 'ARBITRUM' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D39'; Value = '1.203' }
    @{ Cell = 'E39'; Value = '  -0.86%  ' }
    @{ Cell = 'B40'; Value = 'TheSandbox' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = 'D40'; Value = '0.6472' }
    @{ Cell = 'E40'; Value = '  +1.04%  ' }
    @{ Cell = 'B41'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D41'; Value = '4.992' }
    @{ Cell = 'E41'; Value = '  +1.26%  ' }
    @{ Cell = 'B42'; Value = 'TrustWalletToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D42'; Value = '1.222' }
    @{ Cell = 'E42'; Value = '  -5.56%  ' }
    @{ Cell = 'B43'; Value = 'Aptos' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D43'; Value = '11.19' }
    @{ Cell = 'E43'; Value = '  -2.83%  ' }
    @{ Cell = 'B44'; Value = 'Decentraland' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D44'; Value = '0.6124' }
    @{ Cell = 'E44'; Value = '  +1.49%  ' }
    @{ Cell = 'B45'; Value = 'EnergySwap' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D45'; Value = '13.12' }
    @{ Cell = 'E45'; Value = '  -1.16%  ' }
    @{ Cell = 'B46'; Value = 'WEMIXTOKEN' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Cell = 'D46'; Value = '1.289' }
    @{ Cell = 'E46'; Value = '  +0.98%  ' }
    @{ Cell = 'B47'; Value = 'PancakeSwap' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D47'; Value = '3.649' }
    @{ Cell = 'E47'; Value = '  -0.42%  ' }
    @{ Cell = 'B48'; Value = 'NEARProtocol' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D48'; Value = '2.012' }
    @{ Cell = 'E48'; Value = '  +0.28%  ' }
    @{ Cell = 'B49'; Value = 'EOS' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' }
    @{ Cell = 'D49'; Value = '1.210' }
    @{ Cell = 'E49'; Value = '  -2.48%  ' }
    @{ Cell = 'B50'; Value = 'Quant' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = 'D50'; Value = '120.24' }
    @{ Cell = 'E50'; Value = '  -0.94%  ' }
    @{ Cell = 'B51'; Value = 'Aave' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D51'; Value = '78.47' }
    @{ Cell = 'E51'; Value = '  -1.07%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
